$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest
# scraped values. D-column cells whose new text could be parsed as a plain
# number are quote-prefixed (leading "'") so Excel stores them as text
# (preserving exact formatting such as trailing zeros / thousands dots),
# matching the original inline-string cell contents.

$ws.Range("D2").Value = "56.779.43"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "2.505.33"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'496.07"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").Value = "'154.05"
$ws.Range("E6").Value = "  +8.51%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "2.525.37"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'5.80"
$ws.Range("E10").Value = "  +5.55%  "
$ws.Range("D11").Value = "'0.0996"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D14").Value = "2.941.52"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "56.944.08"
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "'21.62"
$ws.Range("E16").Value = "  +4.58%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "2.511.56"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("D20").Value = "'10.40"
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").Value = "'325.52"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").Value = "'59.23"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").Value = "'0.165"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "2.614.52"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("D30").Value = "0.0₃0825"
$ws.Range("E30").Value = "  +3.30%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'152.64"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").Value = "'18.50"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "'1.54"
$ws.Range("E34").Value = "  +3.19%  "
$ws.Range("D35").Value = "'5.29"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("E36").Value = "  +4.04%  "
$ws.Range("D37").Value = "'3.83"
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("D38").Value = "'0.879"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E39").Value = "  +5.13%  "
$ws.Range("D40").Value = "'34.34"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "'3.54"
$ws.Range("E41").Value = "  +3.21%  "
$ws.Range("D42").Value = "'0.0569"
$ws.Range("E42").Value = "  +2.68%  "
$ws.Range("D43").Value = "'0.625"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'5.00"
$ws.Range("E45").Value = "  +7.20%  "
$ws.Range("D46").Value = "'268.67"
$ws.Range("E46").Value = "  +5.77%  "
$ws.Range("D47").Value = "'0.0934"
$ws.Range("E47").Value = "  +2.90%  "
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "'18.01"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("D51").Value = "1.913.15"
$ws.Range("E51").Value = "  -3.94%  "
